# Weekly update: insert a new data row for "Puerro" (Vega Central Mapocho de
# Santiago) by duplicating the first data row of the block (row 62) and
# shifting all the following rows down by one, exactly as a user would do by
# copying row 62 and choosing "Insert Copied Cells" above row 63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlShiftDown = -4121
$xlShiftDown = -4121

# Copy the full row 62 (all columns/styles) and insert it as the new row 63,
# pushing rows 63..101 down to 64..102.
$ws.Rows.Item(62).Copy() | Out-Null
$ws.Rows.Item(63).Insert($xlShiftDown) | Out-Null
